# Apply updated cryptos list values (generated from commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.903.57"
$ws.Range("E2").Value = '  +1.35%  '
$ws.Range("D3").Value = "'2.298.84"
$ws.Range("E3").Value = '  +0.50%  '
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = "'313.15"
$ws.Range("E5").Value = '  -1.59%  '
$ws.Range("D6").Value = "'104.58"
$ws.Range("E6").Value = '  +3.20%  '
$ws.Range("E7").Value = '  -0.48%  '
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("E9").Value = '  +0.23%  '
$ws.Range("D10").Value = "'39.26"
$ws.Range("E10").Value = '  +0.20%  '
$ws.Range("E11").Value = '  +0.34%  '
$ws.Range("E12").Value = '  +0.54%  '
$ws.Range("E13").Value = '  +1.73%  '
$ws.Range("E14").Value = '  +3.36%  '
$ws.Range("D15").Value = "'15.19"
$ws.Range("E15").Value = '  +0.15%  '
$ws.Range("D16").Value = "'2.647.64"
$ws.Range("E16").Value = '  +0.46%  '
$ws.Range("D17").Value = "'2.297.56"
$ws.Range("E17").Value = '  +0.05%  '
$ws.Range("D18").Value = "'42.762.24"
$ws.Range("E18").Value = '  +1.13%  '
$ws.Range("D19").Value = "'7.32"
$ws.Range("E19").Value = '  -1.08%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = "'0.0000105"
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").Value = "'13.60"
$ws.Range("E21").Value = '  +4.06%  '
$ws.Range("D22").Value = "'73.44"
$ws.Range("E22").Value = '  +0.92%  '
$ws.Range("D23").Value = "'3.50"
$ws.Range("E23").Value = '  -1.34%  '
$ws.Range("D24").Value = "'265.42"
$ws.Range("E24").Value = '  -0.70%  '
$ws.Range("E25").Value = '  -1.03%  '
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = '  +0.58%  '
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("D28").Value = "'7.13"
$ws.Range("E28").Value = '  +16.77%  '
$ws.Range("D29").Value = "'2.33"
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("D30").Value = "'22.41"
$ws.Range("E30").Value = '  -0.42%  '
$ws.Range("D31").Value = "'36.14"
$ws.Range("E31").Value = '  -3.43%  '
$ws.Range("D32").Value = "'165.30"
$ws.Range("E32").Value = '  -0.44%  '
$ws.Range("E33").Value = '  -0.73%  '
$ws.Range("D34").Value = "'0.130"
$ws.Range("E34").Value = '  -1.90%  '
$ws.Range("D35").Value = "'2.56"
$ws.Range("E35").Value = '  +0.33%  '
$ws.Range("E36").Value = '  -1.27%  '
$ws.Range("E37").Value = '  -1.13%  '
$ws.Range("D38").Value = "'0.0351"
$ws.Range("E38").Value = '  -1.92%  '
$ws.Range("D39").Value = "'3.74"
$ws.Range("E39").Value = '  +2.53%  '
$ws.Range("D40").Value = "'2.76"
$ws.Range("E40").Value = '  +0.79%  '
$ws.Range("D41").Value = "'1.61"
$ws.Range("E41").Value = '  +5.45%  '
$ws.Range("D42").Value = "'101.15"
$ws.Range("E42").Value = '  +8.83%  '
$ws.Range("D43").Value = "'69.41"
$ws.Range("E43").Value = '  +1.26%  '
$ws.Range("D44").Value = "'0.228"
$ws.Range("E44").Value = '  +1.77%  '
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D46").Value = "'12.18"
$ws.Range("E46").Value = '  +2.07%  '
$ws.Range("D47").Value = "'1.750.62"
$ws.Range("E47").Value = '  +9.41%  '
$ws.Range("D48").Value = "'80.02"
$ws.Range("E48").Value = '  +1.50%  '
$ws.Range("E49").Value = '  -2.50%  '
$ws.Range("D50").Value = "'5.21"
$ws.Range("E50").Value = '  -0.29%  '
$ws.Range("D51").Value = "'8.71"
$ws.Range("E51").Value = '  -2.72%  '
